$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = "21:15"
$ws.Cells.Item(2, 3).Value = 21
$ws.Cells.Item(2, 4).Value = 15
$ws.Cells.Item(2, 5).Value = "PM"
$ws.Cells.Item(2, 6).Value = "Nuit"
# Row 3
$ws.Cells.Item(3, 2).Value = "07:30"
$ws.Cells.Item(3, 3).Value = 7
# Row 4
$ws.Cells.Item(4, 2).Value = "23:15"
$ws.Cells.Item(4, 3).Value = 23
$ws.Cells.Item(4, 4).Value = 15
$ws.Cells.Item(4, 5).Value = "PM"
# Row 5
$ws.Cells.Item(5, 2).Value = "14:30"
$ws.Cells.Item(5, 3).Value = 14
$ws.Cells.Item(5, 4).Value = 30
$ws.Cells.Item(5, 6).Value = "Jour"
# Row 6
$ws.Cells.Item(6, 2).Value = "16:00"
$ws.Cells.Item(6, 3).Value = 16
$ws.Cells.Item(6, 4).Value = 0
$ws.Cells.Item(6, 6).Value = "Jour"
# Row 7
$ws.Cells.Item(7, 2).Value = "15:30"
$ws.Cells.Item(7, 3).Value = 15
$ws.Cells.Item(7, 4).Value = 30
$ws.Cells.Item(7, 5).Value = "PM"
$ws.Cells.Item(7, 6).Value = "Jour"
$ws.Cells.Item(7, 7).Value = "Non"
# Row 8
$ws.Cells.Item(8, 2).Value = "04:00"
$ws.Cells.Item(8, 3).Value = 4
$ws.Cells.Item(8, 4).Value = 0
$ws.Cells.Item(8, 5).Value = "AM"
$ws.Cells.Item(8, 6).Value = "Nuit"
$ws.Cells.Item(8, 7).Value = "Oui"
# Row 9
$ws.Cells.Item(9, 2).Value = "05:15"
$ws.Cells.Item(9, 3).Value = 5
$ws.Cells.Item(9, 4).Value = 15
$ws.Cells.Item(9, 5).Value = "AM"
$ws.Cells.Item(9, 6).Value = "Nuit"
$ws.Cells.Item(9, 7).Value = "Oui"
# Row 10
$ws.Cells.Item(10, 2).Value = "02:45"
$ws.Cells.Item(10, 3).Value = 2
$ws.Cells.Item(10, 4).Value = 45
# Row 11
$ws.Cells.Item(11, 2).Value = "03:00"
$ws.Cells.Item(11, 3).Value = 3
$ws.Cells.Item(11, 4).Value = 0
$ws.Cells.Item(11, 5).Value = "AM"
$ws.Cells.Item(11, 7).Value = "Oui"
# Row 12
$ws.Cells.Item(12, 2).Value = "18:00"
$ws.Cells.Item(12, 3).Value = 18
$ws.Cells.Item(12, 4).Value = 0
# Row 13
$ws.Cells.Item(13, 2).Value = "12:15"
$ws.Cells.Item(13, 3).Value = 12
$ws.Cells.Item(13, 4).Value = 15
# Row 14
$ws.Cells.Item(14, 2).Value = "21:30"
$ws.Cells.Item(14, 3).Value = 21
$ws.Cells.Item(14, 4).Value = 30
$ws.Cells.Item(14, 6).Value = "Nuit"
# Row 15
$ws.Cells.Item(15, 2).Value = "00:45"
$ws.Cells.Item(15, 3).Value = 0
$ws.Cells.Item(15, 6).Value = "Nuit"
$ws.Cells.Item(15, 7).Value = "Oui"
# Row 16
$ws.Cells.Item(16, 2).Value = "00:30"
$ws.Cells.Item(16, 3).Value = 0
$ws.Cells.Item(16, 4).Value = 30
$ws.Cells.Item(16, 5).Value = "AM"
$ws.Cells.Item(16, 6).Value = "Nuit"
$ws.Cells.Item(16, 7).Value = "Oui"
# Row 17
$ws.Cells.Item(17, 2).Value = "03:45"
$ws.Cells.Item(17, 3).Value = 3
$ws.Cells.Item(17, 4).Value = 45
$ws.Cells.Item(17, 5).Value = "AM"
$ws.Cells.Item(17, 6).Value = "Nuit"
$ws.Cells.Item(17, 7).Value = "Oui"
# Row 18
$ws.Cells.Item(18, 2).Value = "03:30"
$ws.Cells.Item(18, 3).Value = 3
$ws.Cells.Item(18, 5).Value = "AM"
$ws.Cells.Item(18, 7).Value = "Oui"
# Row 19
$ws.Cells.Item(19, 2).Value = "18:00"
$ws.Cells.Item(19, 4).Value = 0
# Row 20
$ws.Cells.Item(20, 2).Value = "09:30"
$ws.Cells.Item(20, 3).Value = 9
$ws.Cells.Item(20, 4).Value = 30
$ws.Cells.Item(20, 5).Value = "AM"
# Row 21
$ws.Cells.Item(21, 2).Value = "20:30"
$ws.Cells.Item(21, 3).Value = 20
$ws.Cells.Item(21, 5).Value = "PM"
$ws.Cells.Item(21, 6).Value = "Nuit"
# Row 22
$ws.Cells.Item(22, 2).Value = "11:30"
$ws.Cells.Item(22, 3).Value = 11
$ws.Cells.Item(22, 4).Value = 30
$ws.Cells.Item(22, 5).Value = "AM"
$ws.Cells.Item(22, 6).Value = "Jour"
# Row 23
$ws.Cells.Item(23, 2).Value = "06:15"
$ws.Cells.Item(23, 3).Value = 6
$ws.Cells.Item(23, 4).Value = 15
$ws.Cells.Item(23, 5).Value = "AM"
# Row 24
$ws.Cells.Item(24, 2).Value = "05:30"
$ws.Cells.Item(24, 3).Value = 5
$ws.Cells.Item(24, 4).Value = 30
$ws.Cells.Item(24, 5).Value = "AM"
$ws.Cells.Item(24, 7).Value = "Oui"
# Row 25
$ws.Cells.Item(25, 2).Value = "20:30"
$ws.Cells.Item(25, 3).Value = 20
$ws.Cells.Item(25, 4).Value = 30
$ws.Cells.Item(25, 6).Value = "Nuit"
# Row 26
$ws.Cells.Item(26, 2).Value = "20:30"
$ws.Cells.Item(26, 3).Value = 20
$ws.Cells.Item(26, 4).Value = 30
$ws.Cells.Item(26, 6).Value = "Nuit"
# Row 27
$ws.Cells.Item(27, 2).Value = "19:00"
$ws.Cells.Item(27, 4).Value = 0
# Row 28
$ws.Cells.Item(28, 2).Value = "00:15"
$ws.Cells.Item(28, 3).Value = 0
$ws.Cells.Item(28, 4).Value = 15
$ws.Cells.Item(28, 5).Value = "AM"
$ws.Cells.Item(28, 6).Value = "Nuit"
$ws.Cells.Item(28, 7).Value = "Oui"
# Row 29
$ws.Cells.Item(29, 2).Value = "06:00"
$ws.Cells.Item(29, 3).Value = 6
$ws.Cells.Item(29, 4).Value = 0
$ws.Cells.Item(29, 5).Value = "AM"
$ws.Cells.Item(29, 6).Value = "Jour"
# Row 30
$ws.Cells.Item(30, 2).Value = "09:00"
$ws.Cells.Item(30, 3).Value = 9
$ws.Cells.Item(30, 5).Value = "AM"
# Row 31
$ws.Cells.Item(31, 2).Value = "14:45"
$ws.Cells.Item(31, 3).Value = 14
$ws.Cells.Item(31, 4).Value = 45
$ws.Cells.Item(31, 5).Value = "PM"
$ws.Cells.Item(31, 6).Value = "Jour"
$ws.Cells.Item(31, 7).Value = "Non"
# Row 32
$ws.Cells.Item(32, 2).Value = "13:15"
$ws.Cells.Item(32, 3).Value = 13
$ws.Cells.Item(32, 4).Value = 15
# Row 33
$ws.Cells.Item(33, 2).Value = "07:15"
$ws.Cells.Item(33, 3).Value = 7
$ws.Cells.Item(33, 4).Value = 15
$ws.Cells.Item(33, 6).Value = "Jour"
$ws.Cells.Item(33, 7).Value = "Non"
# Row 34
$ws.Cells.Item(34, 2).Value = "01:15"
$ws.Cells.Item(34, 3).Value = 1
$ws.Cells.Item(34, 4).Value = 15
$ws.Cells.Item(34, 5).Value = "AM"
$ws.Cells.Item(34, 7).Value = "Oui"
# Row 35
$ws.Cells.Item(35, 2).Value = "20:00"
$ws.Cells.Item(35, 3).Value = 20
$ws.Cells.Item(35, 4).Value = 0
$ws.Cells.Item(35, 6).Value = "Nuit"
# Row 36
$ws.Cells.Item(36, 2).Value = "12:45"
$ws.Cells.Item(36, 3).Value = 12
$ws.Cells.Item(36, 4).Value = 45
$ws.Cells.Item(36, 5).Value = "PM"
$ws.Cells.Item(36, 6).Value = "Jour"
$ws.Cells.Item(36, 7).Value = "Non"
# Row 37
$ws.Cells.Item(37, 2).Value = "20:45"
$ws.Cells.Item(37, 3).Value = 20
$ws.Cells.Item(37, 4).Value = 45
$ws.Cells.Item(37, 6).Value = "Nuit"
# Row 38
$ws.Cells.Item(38, 2).Value = "07:00"
$ws.Cells.Item(38, 4).Value = 0
# Row 39
$ws.Cells.Item(39, 2).Value = "18:15"
$ws.Cells.Item(39, 3).Value = 18
$ws.Cells.Item(39, 4).Value = 15
$ws.Cells.Item(39, 5).Value = "PM"
$ws.Cells.Item(39, 6).Value = "Jour"
$ws.Cells.Item(39, 7).Value = "Non"
# Row 40
$ws.Cells.Item(40, 2).Value = "00:30"
$ws.Cells.Item(40, 3).Value = 0
$ws.Cells.Item(40, 4).Value = 30
$ws.Cells.Item(40, 6).Value = "Nuit"
$ws.Cells.Item(40, 7).Value = "Oui"
# Row 41
$ws.Cells.Item(41, 2).Value = "14:30"
$ws.Cells.Item(41, 3).Value = 14
$ws.Cells.Item(41, 4).Value = 30
# Row 42
$ws.Cells.Item(42, 2).Value = "21:30"
$ws.Cells.Item(42, 3).Value = 21
$ws.Cells.Item(42, 5).Value = "PM"
$ws.Cells.Item(42, 7).Value = "Non"
# Row 43
$ws.Cells.Item(43, 2).Value = "16:15"
$ws.Cells.Item(43, 3).Value = 16
$ws.Cells.Item(43, 4).Value = 15
# Row 44
$ws.Cells.Item(44, 2).Value = "20:45"
$ws.Cells.Item(44, 3).Value = 20
$ws.Cells.Item(44, 4).Value = 45
$ws.Cells.Item(44, 6).Value = "Nuit"
# Row 45
$ws.Cells.Item(45, 2).Value = "12:15"
$ws.Cells.Item(45, 3).Value = 12
$ws.Cells.Item(45, 6).Value = "Jour"
# Row 46
$ws.Cells.Item(46, 2).Value = "07:45"
$ws.Cells.Item(46, 3).Value = 7
$ws.Cells.Item(46, 4).Value = 45
$ws.Cells.Item(46, 5).Value = "AM"
$ws.Cells.Item(46, 6).Value = "Jour"
$ws.Cells.Item(46, 7).Value = "Non"
# Row 47
$ws.Cells.Item(47, 2).Value = "18:45"
$ws.Cells.Item(47, 3).Value = 18
$ws.Cells.Item(47, 4).Value = 45
$ws.Cells.Item(47, 5).Value = "PM"
$ws.Cells.Item(47, 6).Value = "Jour"
$ws.Cells.Item(47, 7).Value = "Non"
# Row 48
$ws.Cells.Item(48, 2).Value = "22:45"
$ws.Cells.Item(48, 3).Value = 22
$ws.Cells.Item(48, 5).Value = "PM"
$ws.Cells.Item(48, 7).Value = "Non"
# Row 49
$ws.Cells.Item(49, 2).Value = "22:45"
$ws.Cells.Item(49, 3).Value = 22
$ws.Cells.Item(49, 6).Value = "Nuit"
